$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.104.02"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.927.48"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.64"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.24"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.97"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.70%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.440"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.72"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.412.11"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.992.52"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.72"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.930.28"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "437.22"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.43"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.61"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.02%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.20"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.88"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.60"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.61%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.64"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.53%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.00"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.123"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "42.10"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.00%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "375.56"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.89%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.690.91"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "132.96"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.97"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.00"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.69%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.76%  "
